# Continue the DSA revision sheet: add rows 37-47 (recursion + linked list
# topics, through "binary to integer") and move the active selection to the
# newly added last cell, matching the author's "commiting till linked list
# 3rd question of binary to int" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37: recursion (section heading) ---
$ws.Range("A37").Value = "recursion"

# --- Row 38: power of 2 ---
$ws.Range("A38").Value = "power of 2 "
$ws.Range("B38").Value = "2 se modulo bhi hona chahiye aur uske baad wala bhi 2 se modulo hona chahiye kyunki 2*2*2….. Asia hi hoga and base case true if  n==1 and false  when n == 0;"

# --- Row 39: pow (x,n) ---
$ws.Range("A39").Value = "pow (x,n)"
$ws.Range("B39").Value = "dekho yaar agar n<0 h toh pehle x = 1/x , kardo and n ka abs kardo so who positive jojayega now agar n odd h toh x*fun(x*x,n/2) and even h toh fun(x*x,n/2);  aur base case toh n==0 || x==1 return 1 and x=0 toh 0"
$ws.Range("C39").Value = "imp"

# --- Row 40: reverse a linked list ---
$ws.Range("A40").Value = "reverse a linked list"
$ws.Range("B40").Value = "prev = null , curr = null agar head != NULL h toh curr = head->next; then jab tak curr null nahi hojata loop chalao aur recursion mai head->next = null, haed->next->next = head; aur first mai anshead = reverse(h->n)"
$ws.Range("C40").Value = "imp"

# --- Row 41: merge two sorted ---
$ws.Range("A41").Value = "merge two sorted "
$ws.Range("B41").Value = "ab recursion se kasie karna h nahi pata"
$ws.Range("C41").Value = "imp"

# --- Row 42: reverse a string ---
$ws.Range("A42").Value = "reverse a string"
$ws.Range("B42").Value = "kuch nai wohi h yaar ki bas base case dekhna hota h i>=j h toh return aur swap and I+1 and j-1  that's  it."
$ws.Range("C42").Value = "good"

# Row 43 intentionally left blank (section break, mirrors the gaps already
# used between other topic blocks in this sheet).

# --- Row 44: Linked list (section heading) ---
$ws.Range("A44").Value = "Linked list "

# --- Row 45: delete node ---
$ws.Range("A45").Value = "delete node"
$ws.Range("B45").Value = "kuch nahi h value change kardo aur next wale ko delete kardo ek hi baat h"

# --- Row 46: find middle ---
$ws.Range("A46").Value = "find middle"
$ws.Range("B46").Value = "bruteforce ( size count karo aur fir bas ) = slow fast dono barabar hi h yaar kuch khas diff yahi h ki isme likha kam ha aisa "

# --- Row 47: binary to integer ---
$ws.Range("A47").Value = "binary to integer"
$ws.Range("B47").Value = "bruteforce size nikalo and size - 1 carry hojayega agar 1 aaya toh 2 ki power carry + ans = ans; aur ek aur way h ki ans = ans * 2 + head->val;"
$ws.Range("C47").Value = "imp"

# Scroll/select to match the author's last editing position.
$ws.Range("C47").Select()
